$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Objetivos:" answer text (row 10, columns B and C)
$newObjetivos = "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."
$ws.Range("B10").Value = $newObjetivos
$ws.Range("C10").Value = $newObjetivos

# 2. Insert a new row after row 12 ("Docentes responsáveis:") to hold the professor's
#    name on its own line (blank A, name repeated in B/C).
$ws.Rows("13:13").Insert()

# Copy the formatting from the row above (Objectives row) so the new row gets the
# same non-bold, wrapped B/C styling without inventing new style entries.
$ws.Range("B11").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C11").Copy()
$ws.Range("C13").PasteSpecial(-4122)

$profName = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("B13").Value = $profName
$ws.Range("C13").Value = $profName
$ws.Range("A13").Clear()

# 3. Update "Programa resumido:" answer text (now row 14, columns B and C)
$newResumido = "A definir, de acordo com o tópico programado."
$ws.Range("B14").Value = $newResumido
$ws.Range("C14").Value = $newResumido

# 4. Update "Programa:" answer text (now row 16, columns B and C)
$newPrograma = "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."
$ws.Range("B16").Value = $newPrograma
$ws.Range("C16").Value = $newPrograma

# 5. Update "Bibliografia:" answer text (now row 22 after the earlier row insert shift)
$bibText = "Apostila ou texto fornecido pelo docente responsável. Artigos extraídos de revistas especializadas nas áreas de Ciências e Tecnologia."
$ws.Range("B22").Value = $bibText
$ws.Range("C22").Value = $bibText
